# Insert a new row at position 27 (pushes existing rows 27..61 down to 28..62)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("27:27").Insert()

$ws.Range("A27").Value = 3
$ws.Range("B27").Value = "Femacal de La Calera"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44467
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 100112026
$ws.Range("G27").Value = "Haba"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 14000
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = 14000
$ws.Range("N27").Value = "$/saco 25 kilos"
$ws.Range("O27").Value = "Provincia de Limarí"
$ws.Range("P27").Value = 560
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = "Hortaliza"
